## Inserts a new weekly price-report row for "Bruselas (repollito)" at
## Mercado Mayorista Lo Valledor de Santiago.
##
## The sheet is an append-only daily log: the newest record goes on top
## (row 11, right after the 9 "pinned" header/summary rows), pushing every
## existing data row down by one and growing the used range from
## A1:R60 to A1:R61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 11..60 down to 12..61 (Excel carries the row-11 formatting,
# including the date-number style on column D, into the freshly
# inserted blank row).
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with this week's record.
$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 45022
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 100112035
$ws.Range("G11").Value = "Bruselas (repollito)"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 330
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 22000
$ws.Range("M11").Value = 20970
$ws.Range("N11").Value = '$/malla 20 kilos'
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 1048
$ws.Range("Q11").Value = 20
$ws.Range("R11").Value = "Hortaliza"
